$d = $word.ActiveDocument

# The document contains several "<id>...</id>" markers that were split
# across three separate runs (one run for "<id>", one for the bare
# identifier text, one for "</id>"). The edit collapses each of these
# three runs into a single run per marker, keeping the formatting of the
# opening "<id>" run. The "<id>fig_p064v_1</id>" marker is left untouched.

$targets = @("p064v_1", "p064v_2", "p064v_3", "p064v_4", "p064v_5")

foreach ($id in $targets) {
    $needle = "<id>" + $id + "</id>"
    $rng = $d.Content.Duplicate
    $found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, `
                                $true, 1, $false, $needle, 2)
}
